# "Using Of Where Condition"
# Append 4 new slides to the end of the deck, each a copy of the last
# existing slide (slide 11 - a blank "." placeholder textbox slide),
# matching the pattern already used throughout this deck for new
# lesson slides.

$p = $ppt.ActivePresentation

for ($n = 0; $n -lt 4; $n++) {
    $last = $p.Slides.Item($p.Slides.Count)
    $last.Duplicate() | Out-Null
}
